# Reorder the comma-separated IA Control identifiers in column A for the
# specified rows on the (single, active) worksheet. Only the text content
# of these cells changes - no other formatting/structure is affected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 'AU-4,AU-4 (1)'
$ws.Cells.Item(3, 1).Value = 'SC-5,SC-5 (2),CM-6 b'
$ws.Cells.Item(4, 1).Value = 'CM-5 (1),AC-6 (8),AU-7 b,AU-8 b,AU-12 (3),AU-7 a,AC-6 (9)'
$ws.Cells.Item(5, 1).Value = 'CM-7 b,AC-17 (9),AC-17 (1),CM-6 b'
$ws.Cells.Item(10, 1).Value = 'CM-7 (5) (b),CM-7 (2)'
$ws.Cells.Item(12, 1).Value = 'AC-7 a,AC-7 b'
$ws.Cells.Item(15, 1).Value = 'AU-3 (1),IA-2,IA-8'
$ws.Cells.Item(17, 1).Value = 'MA-4 (1) (a),AU-3,AU-3 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(19, 1).Value = 'IA-5 (1) (a),IA-5 (1) (b),CM-6 b'
$ws.Cells.Item(21, 1).Value = 'MA-4 e,AC-12,MA-4 (7),SC-10'
$ws.Cells.Item(22, 1).Value = 'AU-6 (4),AU-12 a,MA-4 (1) (a),AU-7 (1),AU-3,CM-5 (1),AU-3 (1),CM-6 b,AU-14 (1),AU-7 a'
$ws.Cells.Item(25, 1).Value = 'MA-4 (1) (a),AU-3,AU-3 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(31, 1).Value = 'MA-4 (1) (a),AU-3,AC-2 (4),AU-3 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(39, 1).Value = 'AU-3,CM-6 b'
$ws.Cells.Item(45, 1).Value = 'AC-8 b,AC-8 a,AC-8 c 1, AC-8 c 2, AC-8 c 3'
$ws.Cells.Item(53, 1).Value = 'MA-4 (6),SC-13'
$ws.Cells.Item(56, 1).Value = 'MA-4 (1) (a),AU-12 c'
$ws.Cells.Item(65, 1).Value = 'IA-2 (2),CM-6 b'
$ws.Cells.Item(67, 1).Value = 'MA-4 (1) (a),AU-3,AU-3 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(69, 1).Value = 'AU-12 a,CM-5 (1),AU-7 b,AU-8 b,AU-12 (3),CM-6 b,AU-7 a,AU-12 c'
$ws.Cells.Item(71, 1).Value = 'AU-4 (1),AU-3'
$ws.Cells.Item(77, 1).Value = 'MA-4 (1) (a),AU-3,AC-2 (4),AU-3 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(80, 1).Value = 'IA-2 (1),IA-2 (2),IA-2 (3),IA-2 (4)'
$ws.Cells.Item(81, 1).Value = 'CM-5 (3),CM-6 b'
$ws.Cells.Item(86, 1).Value = 'MA-4 (1) (a),AU-3,AU-3 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(88, 1).Value = 'AU-12 c,CM-5 (1),AC-6 (9),AC-2 (4)'
$ws.Cells.Item(89, 1).Value = 'IA-2,IA-2 (4),IA-2 (5),IA-2 (2),IA-2 (3)'
$ws.Cells.Item(96, 1).Value = 'SC-8 (1),SC-8,AC-18 (1)'
$ws.Cells.Item(97, 1).Value = 'AU-8 b,AU-8 (1) (a),AU-8 (1) (b)'
$ws.Cells.Item(102, 1).Value = 'MA-4 (1) (a),AU-3,AU-3 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(119, 1).Value = 'MA-4 (1) (a),AU-3,AU-3 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(123, 1).Value = 'CM-7 a,CM-7 b'
$ws.Cells.Item(124, 1).Value = 'MA-4 (1) (a),AU-3,AU-3 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(125, 1).Value = 'CM-7 a,AC-18 (1)'
$ws.Cells.Item(128, 1).Value = 'CM-7 a,IA-5 (1) (c),CM-6 b'
$ws.Cells.Item(148, 1).Value = 'MA-4 (1) (a),AU-3,AU-3 (1),AU-14 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(157, 1).Value = 'MA-4 (1) (a),AU-3,AU-3 (1),AU-12 a,AU-12 c'
$ws.Cells.Item(175, 1).Value = 'CM-7 a,SI-16'
